$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.373.49"
$ws.Range("E2").Value = "  +2.05%  "
$ws.Range("D3").Value = "2.353.91"
$ws.Range("E3").Value = "  +1.57%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "542.79"
$ws.Range("E5").Value = "  +1.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.45"
$ws.Range("E6").Value = "  +2.31%  "
$ws.Range("E7").Value = "  +0.62%  "
$ws.Range("E8").Value = "  +4.94%  "
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.63"
$ws.Range("E10").Value = "  +6.24%  "
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.356"
$ws.Range("E12").Value = "  +3.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.91"
$ws.Range("E13").Value = "  +1.60%  "
$ws.Range("D14").Value = "2.774.55"
$ws.Range("E14").Value = "  +1.37%  "
$ws.Range("D15").Value = "58.323.91"
$ws.Range("E15").Value = "  +1.88%  "
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("D17").Value = "2.353.20"
$ws.Range("E17").Value = "  +0.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.73"
$ws.Range("E18").Value = "  +2.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "333.15"
$ws.Range("E19").Value = "  -2.00%  "
$ws.Range("E20").Value = "  +2.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.73"
$ws.Range("E21").Value = "  -2.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("E23").Value = "  +1.62%  "
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.41"
$ws.Range("E27").Value = "  +5.54%  "
$ws.Range("E28").Value = "  +2.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.20"
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("E30").Value = "  +2.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.14"
$ws.Range("E31").Value = "  +0.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.42"
$ws.Range("E32").Value = "  -0.64%  "
$ws.Range("E33").Value = "  +12.68%  "
$ws.Range("E35").Value = "  +6.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.65%  "
$ws.Range("E37").Value = "  -0.68%  "
$ws.Range("E38").Value = "  +4.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "39.19"
$ws.Range("E39").Value = "  +0.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "143.15"
$ws.Range("E40").Value = "  -3.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.379"
$ws.Range("E41").Value = "  +0.48%  "
$ws.Range("E42").Value = "  +1.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "289.55"
$ws.Range("E43").Value = "  +3.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0943"
$ws.Range("E44").Value = "  +1.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.25"
$ws.Range("E45").Value = "  +3.77%  "
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.564"
$ws.Range("E47").Value = "  +1.19%  "
$ws.Range("E48").Value = "  +1.62%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.56"
$ws.Range("E49").Value = "  +1.09%  "
$ws.Range("B50").Value = "Polygon"
$ws.Range("C50").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.382"
$ws.Range("E50").Value = "  +0.50%  "
$ws.Range("E51").Value = "  +0.43%  "
